$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RunManager")

# Update the "Execute" column value for rows 10 and 11 from "Yes" to "No"
$ws.Range("C10").Value = "No"
$ws.Range("C11").Value = "No"

# Update the active selection to match the saved view state (single cell F11)
$ws.Range("F11").Select()
